$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix existing link text (update the URL values, shown as-is in column B)
$ws.Range("B2").Value = "https://www.youtube.com/1"
$ws.Range("B3").Value = "https://www.youtube.com/watch?v=gFIUxGJHXRk2"

# Add new row: "Báo lỗi" label + hyperlink to SharePoint doc
$ws.Range("A4").Value = "Báo lỗi"
$ws.Hyperlinks.Add($ws.Range("B4"), "https://husteduvn-my.sharepoint.com/:x:/g/personal/hinh_nx204650_sis_hust_edu_vn/EX9vkOmM92xJsgH0WFu1WNQBAF8wYkouriFIwiYqTkbdvQ?e=ijYfty") | Out-Null

# Update selection to match target workbook state
$ws.Range("D10").Select()
